$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.771.43'
$ws.Range('E2').Value = '  +0.19%  '

$ws.Range('D3').Value = '3.793.99'
$ws.Range('E3').Value = '  -0.93%  '

$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.20%  '

$ws.Range('D5').Value = '''434.15'
$ws.Range('E5').Value = '  +5.65%  '

$ws.Range('D6').Value = '''139.89'
$ws.Range('E6').Value = '  +6.71%  '

$ws.Range('E7').Value = '  +1.47%  '

$ws.Range('E8').Value = '  -0.16%  '

$ws.Range('E9').Value = '  +0.85%  '

$ws.Range('E10').Value = '  -9.26%  '

$ws.Range('E11').Value = '  -13.06%  '

$ws.Range('D12').Value = '''42.93'
$ws.Range('E12').Value = '  +5.00%  '

$ws.Range('E13').Value = '  +3.71%  '

$ws.Range('D14').Value = '4.398.56'
$ws.Range('E14').Value = '  -0.75%  '

$ws.Range('D15').Value = '''14.85'
$ws.Range('E15').Value = '  -3.40%  '

$ws.Range('D16').Value = '3.822.36'
$ws.Range('E16').Value = '  +0.70%  '

$ws.Range('E17').Value = '  -0.39%  '

$ws.Range('D18').Value = '''19.94'
$ws.Range('E18').Value = '  +1.98%  '

$ws.Range('E19').Value = '  +7.27%  '

$ws.Range('D20').Value = '66.814.85'
$ws.Range('E20').Value = '  -0.37%  '

$ws.Range('D21').Value = '''411.51'
$ws.Range('E21').Value = '  +0.08%  '

$ws.Range('D22').Value = '''14.81'
$ws.Range('E22').Value = '  +2.03%  '

$ws.Range('D23').Value = '''3.25'
$ws.Range('E23').Value = '  +6.68%  '

$ws.Range('D24').Value = '''85.77'
$ws.Range('E24').Value = '  +0.61%  '

$ws.Range('D25').Value = '''36.99'
$ws.Range('E25').Value = '  +0.68%  '

$ws.Range('E26').Value = '  +7.63%  '

$ws.Range('D27').Value = '''9.94'
$ws.Range('E27').Value = '  +39.27%  '

$ws.Range('D28').Value = '''5.55'
$ws.Range('E28').Value = '  -2.61%  '

$ws.Range('D29').Value = '''9.84'
$ws.Range('E29').Value = '  +3.43%  '

$ws.Range('D30').Value = '''729.71'
$ws.Range('E30').Value = '  +4.72%  '

$ws.Range('D31').Value = '''13.82'
$ws.Range('E31').Value = '  +11.08%  '

$ws.Range('E32').Value = '  +9.91%  '

$ws.Range('E33').Value = '  -1.59%  '

$ws.Range('D34').Value = '''42.58'
$ws.Range('E34').Value = '  +10.48%  '

$ws.Range('E35').Value = '  -0.22%  '

$ws.Range('D36').Value = '''0.153'
$ws.Range('E36').Value = '  +0.83%  '

$ws.Range('D37').Value = '''5.62'
$ws.Range('E37').Value = '  +24.68%  '

$ws.Range('D38').Value = '''56.17'
$ws.Range('E38').Value = '  +2.57%  '

$ws.Range('E39').Value = '  +4.39%  '

$ws.Range('E40').Value = '  +44.36%  '

$ws.Range('D41').Value = '''2.95'
$ws.Range('E41').Value = '  -6.90%  '

$ws.Range('E42').Value = '  +3.90%  '

$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').Value = '''0.338'
$ws.Range('E43').Value = '  +14.82%  '

$ws.Range('B44').Value = 'PEPE'
$ws.Range('C44').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D44').Value = '0.0₃0676'
$ws.Range('E44').Value = '  -15.37%  '

$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = '''1.00'
$ws.Range('E45').Value = '  +0.42%  '

$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = '''3.30'
$ws.Range('E46').Value = '  +4.63%  '

$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = '''2.69'
$ws.Range('E47').Value = '  +5.80%  '

$ws.Range('B48').Value = 'LidoDAOToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D48').Value = '''3.32'
$ws.Range('E48').Value = '  +0.74%  '

$ws.Range('E49').Value = '  +0.60%  '

$ws.Range('D50').Value = '''142.19'
$ws.Range('E50').Value = '  -4.49%  '

$ws.Range('E51').Value = '  +1.25%  '
